$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.505.14'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '1.727.79'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('D6').ClearFormats()
$ws.Range('E7').Value = '  +1.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2664'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06217'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.55%  '
$ws.Range('D10').Value = '1.726.12'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('E12').Value = '  +2.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6157'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.514'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.14'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9999'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '26.513.44'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006927'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.65'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').Value = '1.946.74'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.519'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.933'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.279'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.79'
$ws.Range('D25').ClearFormats()
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').Value = '  +2.18%  '
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.82'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.974'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08013'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.708'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.615'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6345'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.54%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9921'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.77%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9253'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.107'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +11.15%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.421'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '105.12'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -7.93%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.006'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01502'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.573'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.50%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3891'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.21%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.925'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +10.36%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1181'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05336'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.89'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.820'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.265'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.05%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3420'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.54%  '
